$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.495.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.478.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.82%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.59%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.546"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.08%  "

$ws.Range("E8").Value = "  +0.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0779"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.65%  "

$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.864.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.37%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.528.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.785"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.356.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0923"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.70%  "

$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.57%  "

$ws.Range("E28").Value = "  -0.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.09%  "

$ws.Range("B34").Value = "ApeXProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.00%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0749"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.38%  "

$ws.Range("E37").Value = "  -4.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.77%  "

$ws.Range("E39").Value = "  -3.25%  "

$ws.Range("E40").Value = "  -7.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.36%  "

$ws.Range("E42").Value = "  +0.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.971.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0282"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.728.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.24%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.85%  "

$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.176"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.18%  "
